# Adds the Modbus "writeLong" command to the Commands sheet, right after
# the existing "writeWord" row (32bit float -> two registers), inserting a
# new row 22 with the writeLong command name/description and pushing the
# remaining command rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# Insert a new row above the current row 22 ("writeSingle...") so the new
# command sits right after "writeWord" (row 21).
$ws.Rows.Item(22).Insert()

$ws.Range("B22").Value = "writeLong(slaveId,register,value) or writeLong([slaveId,register,value],..,[slaveId,register,value])"
$ws.Range("C22").Value = "write 32bit integer to two 16bit int registers: MODBUS function 16"

# Match the row-height auto-fit that happens when the row is inserted.
$ws.Rows.Item(21).RowHeight = 13.8
$ws.Rows.Item(22).RowHeight = 13.8

# Leave the selection where the author left it after adding the row.
$ws.Range("C22").Select()
